# Apply updated cryptos data per commit "Updated cryptos list on Sun Aug 13 16:33:37 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.346.10'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.847.05'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.34'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6288'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07577'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2913'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.51'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07748'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '1.846.21'
$ws.Range('E12').Value = '  -1.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.011'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6775'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001046'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.04'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('D18').Value = '29.347.37'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '228.75'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.433'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.93'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1395'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.439'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.65'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.409'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.469'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05688'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  -0.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.038'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.17%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.153'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.820'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6949'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.581'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01831'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.95%  '
$ws.Range('D38').Value = '1.240.09'
$ws.Range('E38').Value = '  +1.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.715'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.378'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8996'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9991'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').Value = '2.003.96'
$ws.Range('E43').Value = '  -1.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.41'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.47'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.117'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.040'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000115'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1148'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.673'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.85%  '
